# Generate Report for Handoff
#
# Refreshes the localization-status report: the handoff/handback timestamps
# for the files that just became "Ready for handoff" (or whose handback
# transform failed) advance to the latest run, and the duplicate timestamp
# that 2c75f3f7-... previously carried on its own collapses onto the same
# refreshed value as the rest of the batch.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-20-12 02:20:15"
}

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $zhCnRows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-12 02:20:11"
}

# --- de-de sheet: "Latest Handoff Datetime" column (E) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $deDeRows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-12 02:20:15"
}
